function Find-ParaByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $txt = $p.Range.Text.TrimEnd([char]13)
        if ($txt -eq $text) {
            return $p
        }
    }
    return $null
}

function Find-ParaByStart($doc, $startPos) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Start -eq $startPos) {
            return $p
        }
    }
    return $null
}

$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Split "The qt/qml video player should be used. (Need to search.)" into 3
# runs, wrapping "qml" in a spellStart/spellEnd proofErr pair (the proofing
# engine flagging it as an unrecognised word).
$p1 = Find-ParaByText $d "The qt/qml video player should be used. (Need to search.)"
if ($p1 -eq $null) { throw "Change 1: source paragraph not found" }
$rng1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>The qt/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>qml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> video player should be used. (Need to search.)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# --- Change 3 -------------------------------------------------------------
# Split "For file picker, is it okay to use 3rdparty C/C++ APIs? Or is is
# expected to directly use " into multiple runs around grammar/spelling
# proofErr markers, then delete the following "Create thumbnail" paragraph
# (its content is folded away; nothing from it survives in the new text).
$p3 = Find-ParaByText $d "For file picker, is it okay to use 3rdparty C/C++ APIs? Or is is expected to directly use "
if ($p3 -eq $null) { throw "Change 3: source paragraph not found" }
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End)
$xml3 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For file picker, is it okay to use 3rdparty C/C++ APIs? Or </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>is</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> expected to directly use </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng3.InsertXML($xml3)

$pThumb = Find-ParaByText $d "Create thumbnail"
if ($pThumb -eq $null) { throw "Change 3: Create thumbnail paragraph not found" }
$pThumb.Range.Delete()

# --- Change 2 ---------------------------------------------------------
# Insert two new list paragraphs (same numbered list as "10s forward and
# back buttons.", numId=4) right after it:
#   "Video playback window should resize according to the window"
#   "Thumbnail"
$p2 = Find-ParaByText $d "10s forward and back buttons."
if ($p2 -eq $null) { throw "Change 2: anchor paragraph not found" }

$afterP2 = $p2.Range.Duplicate
$afterP2.Collapse(0)
$afterP2.InsertParagraphAfter()

$newP1 = Find-ParaByStart $d $afterP2.Start
if ($newP1 -eq $null) { throw "Change 2: first new paragraph not found" }
$newP1.Range.Text = "Video playback window should resize according to the window"

$newP1 = Find-ParaByText $d "Video playback window should resize according to the window"
if ($newP1 -eq $null) { throw "Change 2: first new paragraph (re-lookup) not found" }
$afterNewP1 = $newP1.Range.Duplicate
$afterNewP1.Collapse(0)
$afterNewP1.InsertParagraphAfter()

$newP2 = Find-ParaByStart $d $afterNewP1.Start
if ($newP2 -eq $null) { throw "Change 2: second new paragraph not found" }
$newP2.Range.Text = "Thumbnail"

Write-Output "All changes applied."
